$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 2 (current NASIRI row) to make room for the
# ZERNAKH ABDELLAH entries, pushing NASIRI down to row 4 and the totals row to row 5.
$ws.Rows("2:3").Insert()

# Row 2: ZERNAKH ABDELLAH first entry
$ws.Range("A2").Value = "ZERNAKH ABDELLAH"
$ws.Range("B2").Value = "IB19558"
$ws.Range("C2").Value = "'145101211406073828000084"
$ws.Range("D2").Value = "MARRAKECH BENI MELLAL"
$ws.Range("E2").Value = "BP"
$ws.Range("F2").Value = "Direction régionale"
$ws.Range("G2").Value = "009/TTT/AV1"
$ws.Range("H2").Value = "mensuelle"
$ws.Range("I2").Value = 4500
$ws.Range("J2").Value = 350
$ws.Range("K2").Value = 4150

# Row 3: ZERNAKH ABDELLAH second entry
$ws.Range("A3").Value = "ZERNAKH ABDELLAH"
$ws.Range("B3").Value = "IB19558"
$ws.Range("C3").Value = "'145101211406073828000084"
$ws.Range("D3").Value = "MARRAKECH BENI MELLAL"
$ws.Range("E3").Value = "BP"
$ws.Range("F3").Value = "Direction régionale"
$ws.Range("G3").Value = "009/TTT/AV1"
$ws.Range("H3").Value = "mensuelle"
$ws.Range("I3").Value = 3500
$ws.Range("J3").Value = 350
$ws.Range("K3").Value = 3150

# Row 5: totals row (spaces for the text columns, summed numbers)
$ws.Range("A5").Value = " "
$ws.Range("B5").Value = " "
$ws.Range("C5").Value = " "
$ws.Range("D5").Value = " "
$ws.Range("E5").Value = " "
$ws.Range("F5").Value = " "
$ws.Range("G5").Value = " "
$ws.Range("H5").Value = " "
$ws.Range("I5").Value = 16500.01
$ws.Range("J5").Value = 1550.01
$ws.Range("K5").Value = 14950
